$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates derived from the target diff (rows 42-45 and 73-78)
$ws.Range("A42").Value = "130886786"
$ws.Range("B42").Value = "57884"
$ws.Range("E42").Value = "100109"
$ws.Range("F42").Value = "Tretåig hackspett"
$ws.Range("G42").Value = "Picoides tridactylus"
$ws.Range("H42").Value = "(Linnaeus, 1758)"
$ws.Range("Q42").Value = "434196"
$ws.Range("R42").Value = "7052215"
$ws.Range("AC42").Value = "Ringhack"
$ws.Range("A43").Value = "130886798"
$ws.Range("Q43").Value = "434094"
$ws.Range("R43").Value = "7052167"
$ws.Range("AC43").Value = "Ringhack äldre"
$ws.Range("A44").Value = "130886828"
$ws.Range("Q44").Value = "434679"
$ws.Range("R44").Value = "7051828"
$ws.Range("A45").Value = "130886846"
$ws.Range("B45").Value = "91828"
$ws.Range("E45").Value = "5432"
$ws.Range("F45").Value = "Granticka"
$ws.Range("G45").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H45").Value = ""
$ws.Range("Q45").Value = "434095"
$ws.Range("R45").Value = "7052227"
$ws.Range("A73").Value = "130886823"
$ws.Range("B73").Value = "57884"
$ws.Range("D73").Value = "NT"
$ws.Range("E73").Value = "100109"
$ws.Range("F73").Value = "Tretåig hackspett"
$ws.Range("G73").Value = "Picoides tridactylus"
$ws.Range("I73").Value = ""
$ws.Range("M73").Value = ""
$ws.Range("N73").Value = ""
$ws.Range("Q73").Value = "434499"
$ws.Range("R73").Value = "7051916"
$ws.Range("AC73").Value = "Ringhack"
$ws.Range("A74").Value = "130886813"
$ws.Range("Q74").Value = "434112"
$ws.Range("R74").Value = "7052117"
$ws.Range("A75").Value = "130886762"
$ws.Range("Q75").Value = "434867"
$ws.Range("R75").Value = "7051762"
$ws.Range("AC75").Value = "Ringhack äldre"
$ws.Range("A76").Value = "130886821"
$ws.Range("Q76").Value = "434468"
$ws.Range("R76").Value = "7051906"
$ws.Range("A77").Value = "130886789"
$ws.Range("Q77").Value = "434159"
$ws.Range("R77").Value = "7052197"
$ws.Range("A78").Value = "130886832"
$ws.Range("B78").Value = "57988"
$ws.Range("D78").Value = "LC"
$ws.Range("E78").Value = "103031"
$ws.Range("F78").Value = "Lavskrika"
$ws.Range("G78").Value = "Perisoreus infaustus"
$ws.Range("I78").Value = "1"
$ws.Range("M78").Value = "födosökande"
$ws.Range("N78").Value = "observerad"
$ws.Range("Q78").Value = "434123"
$ws.Range("R78").Value = "7052111"
$ws.Range("AC78").Value = ""
